$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of rainfall data (dates are Excel serials, same scheme as existing rows).
# Column B on rows 33-37 uses a new bold style (mirrors the xlsx diff: a fresh
# cellXfs entry with fontId=1 "bold" layered on top of the existing wrap/valign style).
$newRows = @(
    @{ Row = 32; Date = 45535; Rain = 12.5; Bold = $false },
    @{ Row = 33; Date = 45536; Rain = 1.5;  Bold = $true  },
    @{ Row = 34; Date = 45537; Rain = 16;   Bold = $true  },
    @{ Row = 35; Date = 45538; Rain = 2.4;  Bold = $true  },
    @{ Row = 36; Date = 45539; Rain = 8;    Bold = $true  },
    @{ Row = 37; Date = 45540; Rain = 7;    Bold = $true  }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $prevRow = $row - 1

    # Copy the formatting (number format / alignment / font) from the row above,
    # then overwrite with this row's values, so every style byte matches the
    # existing A31/B31 template (and we naturally grow cellXfs when a bold
    # variant is needed, exactly like the reference edit did).
    $ws.Range("A" + $prevRow + ":B" + $prevRow).Copy() | Out-Null
    $ws.Range("A" + $row + ":B" + $row).PasteSpecial(-4122) | Out-Null

    $ws.Range("A" + $row).Value = $r.Date
    $ws.Range("B" + $row).Value = $r.Rain

    if ($r.Bold) {
        $ws.Range("B" + $row).Font.Bold = $true
    }
}

$excel.CutCopyMode = $false

# Move the viewport/selection the same way the author's session ended up:
# scrolled down so row 24 is the top visible row, with E33 selected.
$ws.Range("A24").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 24
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E33").Select() | Out-Null
